$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) counts
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 71
$wsExpo.Range("F4").Value = 96
$wsExpo.Range("F7").Value = 2677
$wsExpo.Range("F10").Value = 117
$wsExpo.Range("F11").Value = 10065
$wsExpo.Range("F13").Value = 257
$wsExpo.Range("F15").Value = 623
$wsExpo.Range("F16").Value = 11757
$wsExpo.Range("F17").Value = 12113
$wsExpo.Range("F19").Value = 95

# Sheet "演出" (performances) - update "想去人数" (F column) count
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 7

# Sheet "全部类型" (all types, combined) - update "想去人数" (F column) counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 71
$wsAll.Range("F4").Value = 96
$wsAll.Range("F7").Value = 2677
$wsAll.Range("F8").Value = 7
$wsAll.Range("F11").Value = 117
$wsAll.Range("F12").Value = 10065
$wsAll.Range("F14").Value = 257
$wsAll.Range("F16").Value = 623
$wsAll.Range("F17").Value = 11757
$wsAll.Range("F18").Value = 12113
$wsAll.Range("F20").Value = 95
